$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-6 from 45233 (2023-11-03)
# to 45243 (2023-11-13), keeping existing numeric/date formatting.
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45243
}
